$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-shuffle the Id / Ost / Nord (A/Q/R) values across rows 18-25 and 27-32.
# (Row 26 is left untouched.) Each triple moves as a unit between rows; the
# value that would have mapped into row 18 (111239082) is instead replaced
# with a brand-new Id (111239084).

$ws.Range("A18").Value = 111239084
$ws.Range("Q18").Value = 621795.5397308235
$ws.Range("R18").Value = 7214150.859849748

$ws.Range("A19").Value = 111239096
$ws.Range("Q19").Value = 621773.326892847
$ws.Range("R19").Value = 7214248.527226545

$ws.Range("A20").Value = 111239087
$ws.Range("Q20").Value = 621536.3988506936
$ws.Range("R20").Value = 7214179.606988239

$ws.Range("A21").Value = 111239098
$ws.Range("Q21").Value = 621736.9238461769
$ws.Range("R21").Value = 7214134.075545141

$ws.Range("A22").Value = 111239095
$ws.Range("Q22").Value = 621715.8270385888
$ws.Range("R22").Value = 7214214.866520428

$ws.Range("A23").Value = 111239099
$ws.Range("Q23").Value = 621729.9502675609
$ws.Range("R23").Value = 7214138.867176525

$ws.Range("A24").Value = 111239094
$ws.Range("Q24").Value = 621681.9517352714
$ws.Range("R24").Value = 7214266.364244876

$ws.Range("A25").Value = 111239089
$ws.Range("Q25").Value = 621542.8630217231
$ws.Range("R25").Value = 7214218.370793003

$ws.Range("A27").Value = 111239091
$ws.Range("Q27").Value = 621609.2709173216
$ws.Range("R27").Value = 7214241.392385839

$ws.Range("A28").Value = 111239085
$ws.Range("Q28").Value = 621798.3442589432
$ws.Range("R28").Value = 7214154.78243159

$ws.Range("A29").Value = 111239097
$ws.Range("Q29").Value = 621717.7557529514
$ws.Range("R29").Value = 7214136.675831676

$ws.Range("A30").Value = 111239088
$ws.Range("Q30").Value = 621537.7971145469
$ws.Range("R30").Value = 7214207.587012939

$ws.Range("A31").Value = 111239093
$ws.Range("Q31").Value = 621629.5775533116
$ws.Range("R31").Value = 7214231.645938496

$ws.Range("A32").Value = 111239100
$ws.Range("Q32").Value = 621721.5087325554
$ws.Range("R32").Value = 7214158.82971553
